$d = $word.ActiveDocument

# --- Change 1: "Following a chapter suspension..." was split across two runs;
#     collapse it back into a single run by replacing the combined text. ---
$find1 = $d.Content.Find
$ok1 = $find1.Execute(
    "Following a chapter suspension, grades of former members rebound closer to their expected levels suggesting that suspensions may also improve academics.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Following a chapter suspension, grades of former members rebound closer to their expected levels suggesting that suspensions may also improve academics.",
    2)
if (-not $ok1) { throw "Change 1 (chapter suspension run merge) failed to find target text" }

# --- Change 2: "(Marie and " / "Zölitz" / " 2017)" was split across three runs
#     (with proofErr spell-check wrappers around "Zölitz"); collapse into one run. ---
$find2 = $d.Content.Find
$ok2 = $find2.Execute(
    "(Marie and Zölitz 2017)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "(Marie and Zölitz 2017)",
    2)
if (-not $ok2) { throw "Change 2 (Marie and Zolitz citation run merge) failed to find target text" }

# --- Change 3: append new literature-review bullet paragraphs (a Zotero citation
#     for the NCAA "March Madness" paper plus its outline notes) right after the
#     "not fraternity members)." paragraph, before the trailing empty paragraph. ---
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$insertRange = $lastPara.Range
$insertRange.Collapse(1)
$xmlPackage = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:fldChar w:fldCharType="begin"/>
      </w:r>
      <w:r>
        <w:instrText xml:space="preserve"> ADDIN ZOTERO_ITEM CSL_CITATION {"citationID":"jpzxDsLr","properties":{"formattedCitation":"(\\uc0\\u8220{}MARCH MADNESS: NCAA TOURNAMENT PARTICIPATION AND COLLEGE ALCOHOL USE - White - 2019 - Contemporary Economic Policy - Wiley Online Library\\uc0\\u8221{} n.d.)","plainCitation":"(“MARCH MADNESS: NCAA TOURNAMENT PARTICIPATION AND COLLEGE ALCOHOL USE - White - 2019 - Contemporary Economic Policy - Wiley Online Library” n.d.)","noteIndex":0},"citationItems":[{"id":629,"uris":["http://zotero.org/users/local/Y2pNJapA/items/QSMPGTE4"],"uri":["http://zotero.org/users/local/Y2pNJapA/items/QSMPGTE4"],"itemData":{"id":629,"type":"webpage","title":"MARCH MADNESS: NCAA TOURNAMENT PARTICIPATION AND COLLEGE ALCOHOL USE - White - 2019 - Contemporary Economic Policy - Wiley Online Library","URL":"https://onlinelibrary.wiley.com/doi/full/10.1111/coep.12425","accessed":{"date-parts":[["2021",11,21]]}}}],"schema":"https://github.com/citation-style-language/schema/raw/master/csl-citation.json"} </w:instrText>
      </w:r>
      <w:r>
        <w:fldChar w:fldCharType="separate"/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t>(“MARCH MADNESS: NCAA TOURNAMENT PARTICIPATION AND COLLEGE ALCOHOL USE - White - 2019 - Contemporary Economic Policy - Wiley Online Library” n.d.)</w:t>
      </w:r>
      <w:r>
        <w:fldChar w:fldCharType="end"/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Economic</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Main Point:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">NCAA basketball tournament causes higher levels of drinking for college students. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Increases in drunk driving</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="2"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Data Used: self-reported alcohol information. Harvard Public school CAS data. Bad </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>bad</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>bad</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($xmlPackage)
